$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted as row 250 (Femacal de La Calera - Frutilla),
# pushing the existing rows 250-319 down to 251-320.
$ws.Rows.Item(250).Insert()

# Populate the newly inserted row 250 with the new record's data.
$ws.Range("A250").Value = 3
$ws.Range("B250").Value = "Femacal de La Calera"
$ws.Range("C250").Value = "Coquimbo"
$ws.Range("D250").Value = 44722
$ws.Range("E250").Value = 5
$ws.Range("F250").Value = "Fruta"
$ws.Range("G250").Value = 100101
$ws.Range("H250").Value = "Berries"
$ws.Range("I250").Value = 100112025
$ws.Range("J250").Value = "Frutilla"
$ws.Range("K250").Value = "Sin especificar"
$ws.Range("L250").Value = "Especial"
$ws.Range("M250").Value = 40
$ws.Range("N250").Value = 10000
$ws.Range("O250").Value = 10000
$ws.Range("P250").Value = 10000
$ws.Range("Q250").Value = "$/bandeja 7 kilos"
$ws.Range("R250").Value = "Provincia de Melipilla"
$ws.Range("S250").Value = 1429
$ws.Range("T250").Value = 7
